$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (only B2 changes)
$ws.Range("B2").Value = "NSE:KPRMILL"

# Update column C values (rows 2-10)
$ws.Range("C2").Value = "NSE:ASHIANA"
$ws.Range("C3").Value = "NSE:CCL"
$ws.Range("C4").Value = "NSE:DIGISPICE"
$ws.Range("C5").Value = "NSE:HARDWYN"
$ws.Range("C6").Value = "NSE:INDOCO"
$ws.Range("C7").Value = "NSE:LASA"
$ws.Range("C8").Value = "NSE:ORIENTBELL"
$ws.Range("C9").Value = "NSE:RAILTEL"
$ws.Range("C10").Value = "NSE:RUSHIL"

# Clear column E values for rows 2-7 (they become empty)
$ws.Range("E2:E7").ClearContents()

# Delete rows 11-17 entirely (shifting cells up, removing them from the sheet)
$ws.Range("A11:F17").EntireRow.Delete()
